# Update TPM-derived NATMI metrics for Gdf9-Acvr2a sheet (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.11624
$ws.Range("H2").Value = 6.34872
$ws.Range("I2").Value = 0.1897594766532197
$ws.Range("J2").Value = 0.1897594766532197
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 34.47026378469333
$ws.Range("R2").Value = 310.23237406224
$ws.Range("S2").Value = 0.04130883054773786
$ws.Range("T2").Value = 0.04130883054773787

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.11624
$ws.Range("H3").Value = 6.34872
$ws.Range("I3").Value = 0.1897594766532197
$ws.Range("J3").Value = 0.1897594766532197
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83272
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 58.43130512426666
$ws.Range("R3").Value = 525.8817461184
$ws.Range("S3").Value = 0.0700235106159334
$ws.Range("T3").Value = 0.07002351061593343

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.11624
$ws.Range("H4").Value = 6.34872
$ws.Range("I4").Value = 0.1897594766532197
$ws.Range("J4").Value = 0.1897594766532197
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 55.585847618
$ws.Range("R4").Value = 500.272628562
$ws.Range("S4").Value = 0.06661354187617129
$ws.Range("T4").Value = 0.06661354187617131

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.11624
$ws.Range("H5").Value = 6.34872
$ws.Range("I5").Value = 0.1897594766532197
$ws.Range("J5").Value = 0.1897594766532197
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 9.857884687093334
$ws.Range("R5").Value = 88.72096218384
$ws.Range("S5").Value = 0.01181359361337714
$ws.Range("T5").Value = 0.01181359361337714

# Row 6
$ws.Range("I6").Value = 0.6160274054778138
$ws.Range("J6").Value = 0.6160274054778138
$ws.Range("M6").Value = 16.28844733333333
$ws.Range("N6").Value = 48.865342
$ws.Range("O6").Value = 0.2176904746803693
$ws.Range("P6").Value = 0.2176904746803693
$ws.Range("Q6").Value = 111.90285481355
$ws.Range("R6").Value = 1007.12569332195
$ws.Range("S6").Value = 0.1341032983145816
$ws.Range("T6").Value = 0.1341032983145816

# Row 7
$ws.Range("I7").Value = 0.6160274054778138
$ws.Range("J7").Value = 0.6160274054778138
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83272
$ws.Range("O7").Value = 0.3690119294748028
$ws.Range("P7").Value = 0.3690119294748029
$ws.Range("S7").Value = 0.2273214615047248
$ws.Range("T7").Value = 0.2273214615047248

# Row 8
$ws.Range("I8").Value = 0.6160274054778138
$ws.Range("J8").Value = 0.6160274054778138
$ws.Range("M8").Value = 26.266325
$ws.Range("N8").Value = 78.798975
$ws.Range("O8").Value = 0.3510419771967738
$ws.Range("P8").Value = 0.3510419771967739
$ws.Range("Q8").Value = 180.451622724375
$ws.Range("R8").Value = 1624.064604519375
$ws.Range("S8").Value = 0.2162514784263304
$ws.Range("T8").Value = 0.2162514784263305

# Row 9
$ws.Range("I9").Value = 0.6160274054778138
$ws.Range("J9").Value = 0.6160274054778138
$ws.Range("M9").Value = 4.658207333333333
$ws.Range("N9").Value = 13.974622
$ws.Range("O9").Value = 0.06225561864805391
$ws.Range("P9").Value = 0.06225561864805392
$ws.Range("Q9").Value = 32.00223374555
$ws.Range("R9").Value = 288.02010370995
$ws.Range("S9").Value = 0.03835116723217685
$ws.Range("T9").Value = 0.03835116723217686

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.793503666666667
$ws.Range("H10").Value = 5.380511
$ws.Range("I10").Value = 0.1608202836929164
$ws.Range("J10").Value = 0.1608202836929164
$ws.Range("M10").Value = 16.28844733333333
$ws.Range("N10").Value = 48.865342
$ws.Range("O10").Value = 0.2176904746803693
$ws.Range("P10").Value = 0.2176904746803693
$ws.Range("Q10").Value = 29.21339001664023
$ws.Range("R10").Value = 262.920510149762
$ws.Range("S10").Value = 0.03500904389534262
$ws.Range("T10").Value = 0.03500904389534262

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.793503666666667
$ws.Range("H11").Value = 5.380511
$ws.Range("I11").Value = 0.1608202836929164
$ws.Range("J11").Value = 0.1608202836929164
$ws.Range("M11").Value = 27.61090666666666
$ws.Range("N11").Value = 82.83272
$ws.Range("O11").Value = 0.3690119294748028
$ws.Range("P11").Value = 0.3690119294748029
$ws.Range("Q11").Value = 49.52026234665778
$ws.Range("R11").Value = 445.68236111992
$ws.Range("S11").Value = 0.05934460318420824
$ws.Range("T11").Value = 0.05934460318420825

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.793503666666667
$ws.Range("H12").Value = 5.380511
$ws.Range("I12").Value = 0.1608202836929164
$ws.Range("J12").Value = 0.1608202836929164
$ws.Range("M12").Value = 26.266325
$ws.Range("N12").Value = 78.798975
$ws.Range("O12").Value = 0.3510419771967738
$ws.Range("P12").Value = 0.3510419771967739
$ws.Range("Q12").Value = 47.10875019735833
$ws.Range("R12").Value = 423.978751776225
$ws.Range("S12").Value = 0.05645467036090745
$ws.Range("T12").Value = 0.05645467036090746

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.793503666666667
$ws.Range("H13").Value = 5.380511
$ws.Range("I13").Value = 0.1608202836929164
$ws.Range("J13").Value = 0.1608202836929164
$ws.Range("M13").Value = 4.658207333333333
$ws.Range("N13").Value = 13.974622
$ws.Range("O13").Value = 0.06225561864805391
$ws.Range("P13").Value = 0.06225561864805392
$ws.Range("Q13").Value = 8.35451193242689
$ws.Range("R13").Value = 75.19060739184201
$ws.Range("S13").Value = 0.01001196625245804
$ws.Range("T13").Value = 0.01001196625245805

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.3724043333333333
$ws.Range("H14").Value = 1.117213
$ws.Range("I14").Value = 0.03339283417605023
$ws.Range("J14").Value = 0.03339283417605023
$ws.Range("M14").Value = 16.28844733333333
$ws.Range("N14").Value = 48.865342
$ws.Range("O14").Value = 0.2176904746803693
$ws.Range("P14").Value = 0.2176904746803693
$ws.Range("Q14").Value = 6.065888370205111
$ws.Range("R14").Value = 54.592995331846
$ws.Range("S14").Value = 0.007269301922707234
$ws.Range("T14").Value = 0.007269301922707234

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.3724043333333333
$ws.Range("H15").Value = 1.117213
$ws.Range("I15").Value = 0.03339283417605023
$ws.Range("J15").Value = 0.03339283417605023
$ws.Range("M15").Value = 27.61090666666666
$ws.Range("N15").Value = 82.83272
$ws.Range("O15").Value = 0.3690119294748028
$ws.Range("P15").Value = 0.3690119294748029
$ws.Range("Q15").Value = 10.28242128992889
$ws.Range("R15").Value = 92.54179160935999
$ws.Range("S15").Value = 0.01232235416993643
$ws.Range("T15").Value = 0.01232235416993643

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.3724043333333333
$ws.Range("H16").Value = 1.117213
$ws.Range("I16").Value = 0.03339283417605023
$ws.Range("J16").Value = 0.03339283417605023
$ws.Range("M16").Value = 26.266325
$ws.Range("N16").Value = 78.798975
$ws.Range("O16").Value = 0.3510419771967738
$ws.Range("P16").Value = 0.3510419771967739
$ws.Range("Q16").Value = 9.781693250741666
$ws.Range("R16").Value = 88.035239256675
$ws.Range("S16").Value = 0.01172228653336467
$ws.Range("T16").Value = 0.01172228653336468

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.3724043333333333
$ws.Range("H17").Value = 1.117213
$ws.Range("I17").Value = 0.03339283417605023
$ws.Range("J17").Value = 0.03339283417605023
$ws.Range("M17").Value = 4.658207333333333
$ws.Range("N17").Value = 13.974622
$ws.Range("O17").Value = 0.06225561864805391
$ws.Range("P17").Value = 0.06225561864805392
$ws.Range("Q17").Value = 1.734736596498444
$ws.Range("R17").Value = 15.612629368486
$ws.Range("S17").Value = 0.002078891550041885
$ws.Range("T17").Value = 0.002078891550041885
